# Apply cryptos.xlsx data refresh per commit: "Updated cryptos list on Thu Sep  7 04:34:52 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.843.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.638.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.877.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0496"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.910"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.133.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.808"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.775.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.39%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.418"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.76%  "
